$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ESC")

# Row 7: best_score 0 -> 2, best_time 999999 -> 15
$ws.Range("B7").Value = 2
$ws.Range("D7").Value = 15

# Row 54: best_score 2 -> 12, best_time 13 -> 98
$ws.Range("B54").Value = 12
$ws.Range("D54").Value = 98

# Row 63: best_score 0 -> 1, best_time 999999 -> 20
$ws.Range("B63").Value = 1
$ws.Range("D63").Value = 20
